# Excel COM-interop script replicating the scheduled-runner market-price refresh.
# For each affected Leve row (identified by sheet + row number) this updates the
# currentAveragePrice / currentAveragePriceNQ/HQ / LevePriceNQ/HQ / LeveProfitNQ/HQ
# columns (H..N) with their newly observed values. A few rows also gain or lose a
# trailing N (or M) cell entirely, which is mirrored with ClearContents()/new Value.

$wb = $excel.ActiveWorkbook

# ---- ALC sheet ----
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 4200.3335
$ws.Range("I2").Value = 449.66666
$ws.Range("K2").Value = 449.66666
$ws.Range("M2").Value = -336.66666

# Row 11
$ws.Range("H11").Value = 60
$ws.Range("I11").Value = 60
$ws.Range("K11").Value = 60
$ws.Range("M11").Value = 80

# Row 18
$ws.Range("H18").Value = 989.3
$ws.Range("I18").Value = 929.7273
$ws.Range("J18").Value = 1062.1111
$ws.Range("K18").Value = 929.7273
$ws.Range("L18").Value = 1062.1111
$ws.Range("M18").Value = -645.7273
$ws.Range("N18").Value = -1630.1111

# Row 64
$ws.Range("H64").Value = 4866.3335
$ws.Range("J64").Value = 5214.143
$ws.Range("L64").Value = 5214.143
$ws.Range("N64").Value = -5710.143

# Row 67
$ws.Range("H67").Value = 4866.3335
$ws.Range("J67").Value = 5214.143
$ws.Range("L67").Value = 5214.143
$ws.Range("N67").Value = -6930.143

# Row 70
$ws.Range("H70").Value = 2227.0667
$ws.Range("I70").Value = 2500.4
$ws.Range("J70").Value = 2090.4
$ws.Range("K70").Value = 7501.200000000001
$ws.Range("L70").Value = 6271.200000000001
$ws.Range("M70").Value = -7231.200000000001
$ws.Range("N70").Value = -6811.200000000001

# Row 73
$ws.Range("H73").Value = 2227.0667
$ws.Range("I73").Value = 2500.4
$ws.Range("J73").Value = 2090.4
$ws.Range("K73").Value = 7501.200000000001
$ws.Range("L73").Value = 6271.200000000001
$ws.Range("M73").Value = -6565.200000000001
$ws.Range("N73").Value = -8143.200000000001

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 112
$ws.Range("H112").Value = 4738
$ws.Range("J112").Value = 4880.25
$ws.Range("L112").Value = 14640.75
$ws.Range("N112").Value = -16856.75

# ---- ARM sheet ----
$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 361.8
$ws.Range("I4").Value = 231.66667
$ws.Range("K4").Value = 231.66667
$ws.Range("M4").Value = -115.66667

# Row 5
$ws.Range("H5").Value = 200
$ws.Range("I5").Value = 200
$ws.Range("K5").Value = 200
$ws.Range("M5").Value = -88

# Row 132
$ws.Range("H132").Value = 869.6
$ws.Range("I132").Value = 837.25
$ws.Range("J132").Value = 999
$ws.Range("K132").Value = 2511.75
$ws.Range("L132").Value = 2997
$ws.Range("M132").Value = 18.25
$ws.Range("N132").Value = -8057

# ---- BSM sheet ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 200
$ws.Range("I4").Value = 200
$ws.Range("K4").Value = 200
$ws.Range("M4").Value = -85

# Row 132
$ws.Range("H132").Value = 79999
$ws.Range("J132").Value = 79999
$ws.Range("L132").Value = 79999
$ws.Range("N132").Value = -90119

# Row 134
$ws.Range("H134").Value = 1266
$ws.Range("I134").Value = 1224.5
$ws.Range("J134").Value = 1349
$ws.Range("K134").Value = 3673.5
$ws.Range("L134").Value = 4047
$ws.Range("M134").Value = -1138.5
$ws.Range("N134").Value = -9117

# ---- CRP sheet ----
$ws = $wb.Worksheets.Item("CRP")
# Row 54
$ws.Range("H54").Value = 39999
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

# Row 134
$ws.Range("H134").Value = 4249.5
$ws.Range("I134").Value = 3500
$ws.Range("J134").Value = 4999
$ws.Range("K134").Value = 10500
$ws.Range("L134").Value = 14997
$ws.Range("M134").Value = -7965
$ws.Range("N134").Value = -20067

# ---- CUL sheet ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 8327.375
$ws.Range("I3").Value = 8327.375
$ws.Range("K3").Value = 24982.125
$ws.Range("M3").Value = -24870.125

# Row 95
$ws.Range("H95").Value = 8892.5
$ws.Range("J95").Value = 8892.5
$ws.Range("L95").Value = 26677.5
$ws.Range("N95").Value = -30795.5

# Row 106
$ws.Range("H106").Value = 16483
$ws.Range("I106").Value = 9450
$ws.Range("J106").Value = 19999.5
$ws.Range("K106").Value = 28350
$ws.Range("L106").Value = 59998.5
$ws.Range("M106").Value = -27404
$ws.Range("N106").Value = -61890.5

# Row 118
$ws.Range("H118").Value = 1970.5883
$ws.Range("I118").Value = 1970.5883
$ws.Range("K118").Value = 5911.7649
$ws.Range("M118").Value = -4668.7649

# Row 123
$ws.Range("H123").Value = 19991
$ws.Range("J123").Value = 19991
$ws.Range("L123").Value = 59973
$ws.Range("N123").Value = -64873

# ---- GSM sheet ----
$ws = $wb.Worksheets.Item("GSM")
# Row 18
$ws.Range("H18").Value = 30676
$ws.Range("I18").Value = 24458
$ws.Range("J18").Value = 40003
$ws.Range("K18").Value = 24458
$ws.Range("L18").Value = 40003
$ws.Range("N18").Value = -40589
$ws.Range("M18").Value = -24165

# Row 33
$ws.Range("H33").Value = 10043000
$ws.Range("J33").Value = 10043000
$ws.Range("L33").Value = 10043000
$ws.Range("N33").Value = -10043504

# Row 98
$ws.Range("H98").Value = 14720.2
$ws.Range("J98").Value = 14720.2
$ws.Range("L98").Value = 14720.2
$ws.Range("N98").Value = -20710.2

# Row 102
$ws.Range("H102").Value = 2992.3333
$ws.Range("I102").Value = 2992.3333
$ws.Range("K102").Value = 2992.3333
$ws.Range("M102").Value = -1370.3333

# Row 122
$ws.Range("H122").Value = 2048
$ws.Range("I122").Value = 2048
$ws.Range("K122").Value = 6144
$ws.Range("M122").Value = -3694

# Row 134
$ws.Range("H134").Value = 23379
$ws.Range("J134").Value = 23379
$ws.Range("L134").Value = 70137
$ws.Range("N134").Value = -75207

# ---- LTW sheet ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 7708.9375
$ws.Range("I7").Value = 4593.6665
$ws.Range("K7").Value = 4593.6665
$ws.Range("M7").Value = -4481.6665

# Row 22
$ws.Range("H22").Value = 2026.875
$ws.Range("I22").Value = 2028.6666
$ws.Range("K22").Value = 2028.6666
$ws.Range("M22").Value = -1733.6666

# Row 27
$ws.Range("H27").Value = 2026.875
$ws.Range("I27").Value = 2028.6666
$ws.Range("K27").Value = 2028.6666
$ws.Range("M27").Value = -1921.6666

# Row 40
$ws.Range("H40").Value = 4492.625
$ws.Range("I40").Value = 4181.6665
$ws.Range("J40").Value = 4679.2
$ws.Range("K40").Value = 4181.6665
$ws.Range("L40").Value = 4679.2
$ws.Range("M40").Value = -4045.6665
$ws.Range("N40").Value = -4951.2

# Row 61
$ws.Range("H61").Value = 2799
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2799
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2799
$ws.Range("N61").Value = -3203
$ws.Range("M61").ClearContents()

# Row 68
$ws.Range("H68").Value = 1502.6666
$ws.Range("I68").Value = 1453.3
$ws.Range("J68").Value = 1749.5
$ws.Range("K68").Value = 1453.3
$ws.Range("L68").Value = 1749.5
$ws.Range("M68").Value = -704.3
$ws.Range("N68").Value = -3247.5

# Row 71
$ws.Range("H71").Value = 1502.6666
$ws.Range("I71").Value = 1453.3
$ws.Range("J71").Value = 1749.5
$ws.Range("K71").Value = 7266.5
$ws.Range("L71").Value = 8747.5
$ws.Range("M71").Value = -3522.5
$ws.Range("N71").Value = -16235.5

# Row 100
$ws.Range("H100").Value = 3609.45
$ws.Range("J100").Value = 3474
$ws.Range("L100").Value = 3474
$ws.Range("N100").Value = -4556

# Row 113
$ws.Range("H113").Value = 2799
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2799
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2799
$ws.Range("N113").Value = -7139
$ws.Range("M113").ClearContents()

# Row 122
$ws.Range("H122").Value = 6634.64
$ws.Range("I122").Value = 5248.25
$ws.Range("K122").Value = 15744.75
$ws.Range("M122").Value = -13294.75

# Row 126
$ws.Range("H126").Value = 7708.9375
$ws.Range("I126").Value = 4593.6665
$ws.Range("K126").Value = 13780.9995
$ws.Range("M126").Value = -11310.9995

# Row 127
$ws.Range("H127").Value = 60000
$ws.Range("J127").Value = 60000
$ws.Range("L127").Value = 60000
$ws.Range("N127").Value = -69920

# ---- WVR sheet ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 2144.625
$ws.Range("I132").Value = 2144.625
$ws.Range("K132").Value = 6433.875
$ws.Range("M132").Value = -3903.875

# Row 136
$ws.Range("H136").Value = 1120.3334
$ws.Range("I136").Value = 1181.6666
$ws.Range("J136").Value = 997.6667
$ws.Range("K136").Value = 3544.9998
$ws.Range("L136").Value = 2993.0001
$ws.Range("M136").Value = -994.9998000000001
$ws.Range("N136").Value = -8093.0001
